$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing "Is Significant" values (column F) before overwriting
$isSig = @()
for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $isSig += $ws.Range("F$row").Value()
}

# New header labels: existing F1 becomes "Observed", add G1 "Expected" and H1 "Is Significant"
$ws.Range("F1").Value = "Observed"
$ws.Range("G1").Value = "Expected"
$ws.Range("H1").Value = "Is Significant"

# Copy header style (bold/border/center) from F1 onto the new G1:H1 header cells
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$observed = @(
    "[707  45] ; [283 105]",
    "[632 122] ; [266 132]",
    "[646 102] ; [295 101]",
    "[641 109] ; [292 100]",
    "[677  77] ; [320  78]",
    "[713  36] ; [347  47]",
    "[714  33] ; [353  41]",
    "[702  42] ; [348  45]",
    "[549 203] ; [250 146]"
)

$expected = @(
    "[653.05263158  98.94736842] ; [336.94736842  51.05263158]",
    "[587.75347222 166.24652778] ; [310.24652778  87.75347222]",
    "[615.26923077 132.73076923] ; [325.73076923  70.26923077]",
    "[612.7408056 137.2591944] ; [320.2591944  71.7408056]",
    "[652.55034722 101.44965278] ; [344.44965278  53.55034722]",
    "[694.61067367  54.38932633] ; [365.38932633  28.61067367]",
    "[698.55302366  48.44697634] ; [368.44697634  25.55302366]",
    "[687.07124011  56.92875989] ; [362.92875989  30.07124011]",
    "[523.38675958 228.61324042] ; [275.61324042 120.38675958]"
)

for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value = $observed[$i]
}
for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $expected[$i]
}
for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws.Range("H$row").Value = $isSig[$i]
}

$ws.Range("A1").Select() | Out-Null
